# "corrected data cleaning for pre/post/total fixation data"
# - Header row (row 1) loses its bold/bordered/centered style -> plain "Normal" style,
#   and the "Unnamed: 0" label in A1 is cleared to blank.
# - Several numeric metrics in rows 3-7 are recalculated (data-cleaning fix).
# - The trailing, entirely-blank row 11 is removed (used range shrinks to A1:R10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: drop direct formatting (bold font + thin border + center/top align) ---
$ws.Range("A1:R1").Style = "Normal"

# --- A1 no longer holds the "Unnamed: 0" label ---
$ws.Range("A1").ClearContents()

# --- Row 3 (Revisit count) ---
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 30
$ws.Range("E3").Value = 4
$ws.Range("G3").Value = 6
$ws.Range("H3").Value = 10
$ws.Range("K3").Value = 16
$ws.Range("L3").Value = 3
$ws.Range("N3").Value = 27
$ws.Range("O3").Value = 2

# --- Row 4 (Fixation count) ---
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 64
$ws.Range("E4").Value = 6
$ws.Range("G4").Value = 9
$ws.Range("H4").Value = 14
$ws.Range("K4").Value = 24
$ws.Range("L4").Value = 7
$ws.Range("N4").Value = 80
$ws.Range("O4").Value = 3

# --- Row 5 (Dwell time (ms)) ---
$ws.Range("B5").Value = 116.8
$ws.Range("C5").Value = 18026.64
$ws.Range("E5").Value = 1201.06
$ws.Range("G5").Value = 1659.91
$ws.Range("H5").Value = 4671.42
$ws.Range("K5").Value = 8617.45
$ws.Range("L5").Value = 3120.12
$ws.Range("N5").Value = 25342.61
$ws.Range("O5").Value = 767.36

# --- Row 6 (Dwell time (%)) ---
$ws.Range("B6").Value = 0.1
$ws.Range("C6").Value = 14.8
$ws.Range("D6").Value = 0.2
$ws.Range("E6").Value = 0.99
$ws.Range("G6").Value = 1.36
$ws.Range("H6").Value = 3.84
$ws.Range("I6").Value = 0.69
$ws.Range("J6").Value = 0.3
$ws.Range("K6").Value = 7.08
$ws.Range("L6").Value = 2.56
$ws.Range("M6").Value = 0.2
$ws.Range("N6").Value = 20.81
$ws.Range("O6").Value = 0.63
$ws.Range("P6").Value = 0.69
$ws.Range("R6").Value = 0.14

# --- Row 7 (Fixation duration (ms)) ---
$ws.Range("B7").Value = 116.8
$ws.Range("C7").Value = 281.67
$ws.Range("E7").Value = 200.18
$ws.Range("G7").Value = 184.43
$ws.Range("H7").Value = 333.67
$ws.Range("K7").Value = 359.06
$ws.Range("L7").Value = 445.73
$ws.Range("N7").Value = 316.78
$ws.Range("O7").Value = 255.79

# --- Drop the trailing blank row (old row 11); dimension becomes A1:R10 ---
$ws.Rows.Item(11).Delete()
